$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

$ws.Range("A7").Value = "projectFolderID"
$ws.Range("B7").Value = "12b1y-sg6E5rox-ntC2UustqiREhdGE5X"
$ws.Range("C7").Value = "ID of project ID from Drive"

$ws.Range("C8").Select()
